# "corregi la escala del flux" -- fix the Flux (column C) unit scale on the
# h3 sheet and add a corrected copy of h2 (named h6) with the same fix,
# mirroring the author's workflow of creating a fixed copy while patching
# the duplicate data already present on h3.

$wb = $excel.ActiveWorkbook

# Corrected Flux values (old value * 0.3213) for rows 2..43, column C.
$newFlux = @(
    11.7340520547945, 18.2920931506849, 31.9099315068493, 23.3800767123288,
    54.4185369863015, 11.9717260273972, 15.6336657534247, 22.6406465753425,
    30.7831808219178, 47.1298684931506, 12.8696054794521, 42.1563205479453,
    32.5965452054796, 16.3995041095891, 38.9961369863015, 99.9551095890411,
    72.3321123287672, 1.13694811321849, 1.09897416078272, 0.158523097967821,
    2.26850205797895, 1.08678176214179, 0.111097662391457, 2.04482995651593,
    1.64335124545848, 2.52638630136986, 0.598586301369864, 2.00702465753425,
    0.431334246575343, 2.42075342465754, 1.93660273972603, 2.43835890410959,
    2.98041470641543, 0.941893150684932, 1.13905246172262, 0.82662503603988,
    0.862668493150683, 0.572178082191781, 0.564672291520741, 0.827457534246574,
    5.40696718926627, 1.00351232876712
)

# --- 1. Fix the Flux column on h3 (rows 2..43) ----------------------------
$h3 = $wb.Worksheets.Item("h3")
for ($i = 0; $i -lt $newFlux.Length; $i++) {
    $row = $i + 2
    $h3.Range("C$row").Value = $newFlux[$i]
}

# --- 2. Create h6 as a corrected copy of h2 --------------------------------
$h2 = $wb.Worksheets.Item("h2")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$null = $h2.Copy($null, $lastSheet)
$h6 = $wb.Worksheets.Item($wb.Worksheets.Count)
$h6.Name = "h6"

for ($i = 0; $i -lt $newFlux.Length; $i++) {
    $row = $i + 2
    $h6.Range("C$row").Value = $newFlux[$i]
}

# --- 3. View / selection bookkeeping ---------------------------------------
$null = $h2.Range("F1").Select()
$null = $h6.Range("H1").Select()

$null = $h3.Activate()
$null = $h3.Range("S9").Select()
